# Update dashboards - 2026-01-01
# Rolling-window update: each indicator's "as-of" date advances by one day
# (and the latest-5 observations window Q:U shifts left, dropping the
# oldest reading and bringing in a new one at Q).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 29 - 5yr, 5yr Forward (T5YIFR)
$ws.Range("N29").Value = 46022
$ws.Range("Q29").Value = 2.24
$ws.Range("R29").Value = 2.23
$ws.Range("S29").Value = 2.21

# Row 30 - 10yr TIPS (T10YIE)
$ws.Range("N30").Value = 46022
$ws.Range("Q30").Value = 2.25
$ws.Range("R30").Value = 2.24
$ws.Range("S30").Value = 2.22
$ws.Range("T30").Value = 2.23

# Row 47 - FFR (DFF)
$ws.Range("N47").Value = 46021

# Row 48 - 2y UST (DGS2)
$ws.Range("N48").Value = 46021
$ws.Range("R48").Value = 3.45
$ws.Range("S48").Value = 3.46
$ws.Range("T48").Value = 3.47
$ws.Range("U48").Value = 3.48

# Row 49 - 5y UST (DGS5)
$ws.Range("N49").Value = 46021
$ws.Range("Q49").Value = 3.68
$ws.Range("R49").Value = 3.67
$ws.Range("S49").Value = 3.68
$ws.Range("T49").Value = 3.7
$ws.Range("U49").Value = 3.72

# Row 50 - 10y UST (DGS10)
$ws.Range("N50").Value = 46021
$ws.Range("Q50").Value = 4.14
$ws.Range("R50").Value = 4.12
$ws.Range("S50").Value = 4.14
$ws.Range("T50").Value = 4.15
$ws.Range("U50").Value = 4.18

# Row 51 - 30y Mtg. (MORTGAGE30US) - weekly series, advances a full week
# and picks up the highlighted-date format used by the other rows (N47:N50,N52).
$ws.Range("N50").Copy()
$ws.Range("N51").PasteSpecial(-4122)
$ws.Range("N51").Value = 46020
$ws.Range("Q51").Value = 6.15
$ws.Range("R51").Value = 6.18
$ws.Range("S51").Value = 6.21
$ws.Range("T51").Value = 6.22
$ws.Range("U51").Value = 6.19

# Row 52 - BAA (DBAA)
$ws.Range("N52").Value = 46021
$ws.Range("Q52").Value = 5.89
$ws.Range("R52").Value = 5.88
$ws.Range("S52").Value = 5.89
$ws.Range("T52").Value = 5.88
$ws.Range("U52").Value = 5.92
